$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the header text in A1 (string value is unchanged, "HK_G_acc_LT"),
# which causes a new shared-string entry to be (re)created for it.
$ws.Range("A1").Value = "HK_G_acc_LT"

# Scale every numeric value in A2:A49 by the factor 555/551 (new sample size),
# matching the recomputed percentages for the new HK genes/reactions set.
$factor = 555.0 / 551.0
for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value = $old * $factor
}
